# saldosporproducto.xlsx - "continuacion scripts saldos por producto"
#
# Adds a new "escenario" column (Q) to the "Datos" sheet and repurposes the
# first data row (row 2) to represent the "Todos los productos" scenario:
#   - K2, L2, N2, O2 (cuentas/tarjetasCredito/inversiones/crediagil ids) are
#     cleared
#   - M2 (creditos id) becomes the scenario code "4676"
#   - Q1 gets the new header "escenario"
#   - Q2 gets the new value "Todos los productos"
#   - Q3:Q14 get the same cell formatting as the neighbouring P column so the
#     new column matches the look of the rest of the table
#   - the view scrolls right a bit and the header row Q1 is (re)selected

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Datos")

# --- new column width -----------------------------------------------------
$ws.Columns.Item(17).ColumnWidth = 18.71

# --- header row -------------------------------------------------------------
$ws.Range("Q1").Value = "escenario"
$ws.Range("P1").Copy()
$ws.Range("Q1").PasteSpecial(-4122)

# --- row 2: "Todos los productos" scenario ---------------------------------
$ws.Range("K2").ClearContents()
$ws.Range("L2").ClearContents()
$ws.Range("M2").Value = "4676"
$ws.Range("N2").ClearContents()
$ws.Range("O2").ClearContents()
$ws.Range("Q2").Value = "Todos los productos"

# match the row-2 cell styling (border/number format/alignment) used by the
# rest of that row
$ws.Range("K2").Copy()
$ws.Range("L2").PasteSpecial(-4122)
$ws.Range("K2").Copy()
$ws.Range("N2").PasteSpecial(-4122)
$ws.Range("K2").Copy()
$ws.Range("O2").PasteSpecial(-4122)
$ws.Range("P2").Copy()
$ws.Range("Q2").PasteSpecial(-4122)

# --- rows 3-14: give column Q the same look as column P in that row --------
$rows = 3..14
foreach ($r in $rows) {
    $ws.Range("P$r").Copy()
    $ws.Range("Q$r").PasteSpecial(-4122)
}

$excel.CutCopyMode = $false

# --- view state --------------------------------------------------------------
$win = $excel.ActiveWindow
$win.ScrollColumn = 13
$ws.Range("A1:Q1").Select()
